# Update BIBI PF annual retention metrics (row 36 and row 37)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 36: num_customers 116 -> 117, retention_rate recalculated (117/1930)
$ws.Range("C36").Value = 117
$ws.Range("E36").Value = 0.06062176165803109

# Row 37: num_customers 720 -> 724, cohort_size 720 -> 724
$ws.Range("C37").Value = 724
$ws.Range("D37").Value = 724
